# Auto-generated edit script applying numeric updates to the Leve profit sheets.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across rows
# in multiple worksheets, matching the upstream scheduled-runner data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3450.5334
$ws.Range("I43").Value = 1580.2
$ws.Range("K43").Value = 1580.2
$ws.Range("M43").Value = -1511.2

$ws.Range("H51").Value = 5885.7144
$ws.Range("J51").Value = 9247.5
$ws.Range("L51").Value = 9247.5
$ws.Range("N51").Value = -10215.5

$ws.Range("H112").Value = 2419
$ws.Range("I112").Value = 656.6667
$ws.Range("J112").Value = 2749.4375
$ws.Range("K112").Value = 1970.0001
$ws.Range("L112").Value = 8248.3125
$ws.Range("M112").Value = -862.0001
$ws.Range("N112").Value = -10464.3125

$ws.Range("H132").Value = 32703.521
$ws.Range("I132").Value = 32703.521
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 98110.56299999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -95580.56299999999
$ws.Range("N132").ClearContents()

$ws.Range("H135").Value = 3638.4546
$ws.Range("I135").Value = 2574.7144
$ws.Range("K135").Value = 23172.4296
$ws.Range("M135").Value = -20637.4296

$ws.Range("H137").Value = 35855.668
$ws.Range("I137").Value = 100701.336
$ws.Range("K137").Value = 302104.008
$ws.Range("M137").Value = -299554.008

$ws.Range("H138").Value = 17586.92
$ws.Range("I138").Value = 1005.125
$ws.Range("J138").Value = 74438.78999999999
$ws.Range("K138").Value = 3015.375
$ws.Range("L138").Value = 223316.37
$ws.Range("M138").Value = 2124.625
$ws.Range("N138").Value = -233596.37

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7932.25
$ws.Range("I61").Value = 1184.6666
$ws.Range("K61").Value = 1184.6666
$ws.Range("M61").Value = -972.6666

$ws.Range("H88").Value = 6433.4546
$ws.Range("I88").Value = 2217.25
$ws.Range("J88").Value = 8842.714
$ws.Range("K88").Value = 2217.25
$ws.Range("L88").Value = 8842.714
$ws.Range("M88").Value = -1811.25
$ws.Range("N88").Value = -9654.714

$ws.Range("H91").Value = 6433.4546
$ws.Range("I91").Value = 2217.25
$ws.Range("J91").Value = 8842.714
$ws.Range("K91").Value = 2217.25
$ws.Range("L91").Value = 8842.714
$ws.Range("M91").Value = -813.25
$ws.Range("N91").Value = -11650.714

$ws.Range("H132").Value = 1119.591
$ws.Range("I132").Value = 884.439
$ws.Range("J132").Value = 4333.3335
$ws.Range("K132").Value = 2653.317
$ws.Range("L132").Value = 13000.0005
$ws.Range("M132").Value = -123.317
$ws.Range("N132").Value = -18060.0005

$ws.Range("H136").Value = 7932.25
$ws.Range("I136").Value = 1184.6666
$ws.Range("K136").Value = 3553.9998
$ws.Range("M136").Value = -1003.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 1683
$ws.Range("I36").Value = 1439.6
$ws.Range("J36").Value = 2900
$ws.Range("K36").Value = 1439.6
$ws.Range("L36").Value = 2900
$ws.Range("M36").Value = -905.5999999999999
$ws.Range("N36").Value = -3968

$ws.Range("H86").Value = 1538.6111
$ws.Range("J86").Value = 2500
$ws.Range("L86").Value = 2500
$ws.Range("N86").Value = -4746

$ws.Range("H89").Value = 1538.6111
$ws.Range("J89").Value = 2500
$ws.Range("L89").Value = 12500
$ws.Range("N89").Value = -23732

$ws.Range("H123").Value = 41000

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 203153
$ws.Range("I6").Value = 3500
$ws.Range("J6").Value = 219790.75
$ws.Range("K6").Value = 3500
$ws.Range("L6").Value = 219790.75
$ws.Range("M6").Value = -3387
$ws.Range("N6").Value = -220016.75

$ws.Range("H16").Value = 2237.4167
$ws.Range("I16").Value = 1426.3334
$ws.Range("J16").Value = 4670.6665
$ws.Range("K16").Value = 1426.3334
$ws.Range("L16").Value = 4670.6665
$ws.Range("M16").Value = -1139.3334
$ws.Range("N16").Value = -5244.6665

$ws.Range("H31").Value = 5001264
$ws.Range("I31").Value = 6250938.5
$ws.Range("J31").Value = 2566
$ws.Range("K31").Value = 6250938.5
$ws.Range("L31").Value = 2566
$ws.Range("M31").Value = -6250643.5
$ws.Range("N31").Value = -3156

$ws.Range("H34").Value = 5001264
$ws.Range("I34").Value = 6250938.5
$ws.Range("J34").Value = 2566
$ws.Range("K34").Value = 6250938.5
$ws.Range("L34").Value = 2566
$ws.Range("M34").Value = -6250736.5
$ws.Range("N34").Value = -2970

$ws.Range("H86").Value = 40782.43
$ws.Range("I86").Value = 61786.082
$ws.Range("J86").Value = 12777.556
$ws.Range("K86").Value = 61786.082
$ws.Range("L86").Value = 12777.556
$ws.Range("M86").Value = -60663.082
$ws.Range("N86").Value = -15023.556

$ws.Range("H89").Value = 40782.43
$ws.Range("I89").Value = 61786.082
$ws.Range("J89").Value = 12777.556
$ws.Range("K89").Value = 308930.41
$ws.Range("L89").Value = 63887.78
$ws.Range("M89").Value = -303314.41
$ws.Range("N89").Value = -75119.78

$ws.Range("H99").Value = 11720.889
$ws.Range("I99").Value = 9621.25
$ws.Range("K99").Value = 9621.25
$ws.Range("M99").Value = -8123.25

$ws.Range("H107").Value = 863.05884
$ws.Range("I107").Value = 787.3
$ws.Range("K107").Value = 787.3
$ws.Range("M107").Value = 1132.7

$ws.Range("H113").Value = 2237.4167
$ws.Range("I113").Value = 1426.3334
$ws.Range("J113").Value = 4670.6665
$ws.Range("K113").Value = 1426.3334
$ws.Range("L113").Value = 4670.6665
$ws.Range("M113").Value = 743.6666
$ws.Range("N113").Value = -9010.666499999999

$ws.Range("H123").Value = 35000
$ws.Range("J123").Value = 35000
$ws.Range("L123").Value = 35000
$ws.Range("N123").Value = -44800

$ws.Range("H126").Value = 11720.889
$ws.Range("I126").Value = 9621.25
$ws.Range("K126").Value = 28863.75
$ws.Range("M126").Value = -26393.75

$ws.Range("H134").Value = 3159.4
$ws.Range("I134").Value = 2988.2354
$ws.Range("K134").Value = 8964.706200000001
$ws.Range("M134").Value = -6429.706200000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 825.2632
$ws.Range("J5").Value = 1215
$ws.Range("L5").Value = 3645
$ws.Range("N5").Value = -3869

$ws.Range("H107").Value = 2117.7896
$ws.Range("I107").Value = 4498.5713
$ws.Range("K107").Value = 13495.7139
$ws.Range("M107").Value = -11575.7139

$ws.Range("H135").Value = 825.2632
$ws.Range("J135").Value = 1215
$ws.Range("L135").Value = 10935
$ws.Range("N135").Value = -16005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 12964.308
$ws.Range("I80").Value = 7417.125
$ws.Range("J80").Value = 21839.8
$ws.Range("K80").Value = 7417.125
$ws.Range("L80").Value = 21839.8
$ws.Range("M80").Value = -6419.125
$ws.Range("N80").Value = -23835.8

$ws.Range("H83").Value = 12964.308
$ws.Range("I83").Value = 7417.125
$ws.Range("J83").Value = 21839.8
$ws.Range("K83").Value = 37085.625
$ws.Range("L83").Value = 109199
$ws.Range("M83").Value = -32093.625
$ws.Range("N83").Value = -119183

$ws.Range("H126").Value = 2198.3333
$ws.Range("I126").Value = 1638
$ws.Range("K126").Value = 4914
$ws.Range("M126").Value = -2444

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1085.2222
$ws.Range("J82").Value = 1633.75
$ws.Range("L82").Value = 1633.75
$ws.Range("N82").Value = -2355.75

$ws.Range("H85").Value = 1085.2222
$ws.Range("J85").Value = 1633.75
$ws.Range("L85").Value = 1633.75
$ws.Range("N85").Value = -4129.75

$ws.Range("H93").Value = 1026.6296
$ws.Range("I93").Value = 720.5217
$ws.Range("J93").Value = 2786.75
$ws.Range("K93").Value = 720.5217
$ws.Range("L93").Value = 2786.75
$ws.Range("M93").Value = 527.4783
$ws.Range("N93").Value = -5282.75

$ws.Range("H132").Value = 3446.9412
$ws.Range("I132").Value = 3372.8667
$ws.Range("J132").Value = 4002.5
$ws.Range("K132").Value = 10118.6001
$ws.Range("L132").Value = 12007.5
$ws.Range("M132").Value = -7588.6001
$ws.Range("N132").Value = -17067.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 430.17392
$ws.Range("I100").Value = 420.66666
$ws.Range("J100").Value = 530
$ws.Range("K100").Value = 841.33332
$ws.Range("L100").Value = 1060
$ws.Range("M100").Value = -300.33332
$ws.Range("N100").Value = -2142

$ws.Range("H115").Value = 35000
$ws.Range("J115").Value = 35000
$ws.Range("L115").Value = 35000
$ws.Range("N115").Value = -38134

